$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.627.19"
$ws.Range("E2").Value = "  -1.18%  "
$ws.Range("D3").Value = "1.617.30"
$ws.Range("E3").Value = "  -1.15%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.991"
$ws.Range("E4").Value = "  -0.92%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "209.38"
$ws.Range("E5").Value = "  -1.37%  "
$ws.Range("E6").Value = "  -1.50%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.990"
$ws.Range("E7").Value = "  -0.96%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.16"
$ws.Range("E8").Value = "  -0.69%  "
$ws.Range("E9").Value = "  -1.13%  "
$ws.Range("E10").Value = "  -1.58%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0874"
$ws.Range("E11").Value = "  -1.04%  "
$ws.Range("D12").Value = "1.847.36"
$ws.Range("E12").Value = "  -1.11%  "
$ws.Range("D13").Value = "1.627.48"
$ws.Range("E13").Value = "  -0.49%  "
$ws.Range("E14").Value = "  -1.91%  "
$ws.Range("E15").Value = "  -1.68%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.62"
$ws.Range("E16").Value = "  -1.22%  "
$ws.Range("D17").Value = "27.650.68"
$ws.Range("E17").Value = "  -1.12%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "227.33"
$ws.Range("E18").Value = "  -1.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.66"
$ws.Range("E19").Value = "  +1.34%  "
$ws.Range("D20").Value = "0.0₃0715"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.991"
$ws.Range("E21").Value = "  -0.88%  "
$ws.Range("E22").Value = "  -1.51%  "
$ws.Range("E23").Value = "  -3.03%  "
$ws.Range("E24").Value = "  -1.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.32"
$ws.Range("E25").Value = "  -0.37%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.88"
$ws.Range("E26").Value = "  -1.25%  "
$ws.Range("E27").Value = "  -1.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.42"
$ws.Range("E28").Value = "  -1.60%  "
$ws.Range("E29").Value = "  -0.90%  "
$ws.Range("E30").Value = "  -1.28%  "
$ws.Range("E31").Value = "  -1.22%  "
$ws.Range("E32").Value = "  -1.00%  "
$ws.Range("E33").Value = "  -0.21%  "
$ws.Range("D34").Value = "1.390.69"
$ws.Range("E34").Value = "  -1.22%  "
$ws.Range("E35").Value = "  +1.46%  "
$ws.Range("E36").Value = "  -1.19%  "
$ws.Range("E37").Value = "  -1.60%  "
$ws.Range("E38").Value = "  -0.10%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.554"
$ws.Range("E39").Value = "  -1.64%  "
$ws.Range("E40").Value = "  -3.43%  "
$ws.Range("E41").Value = "  -1.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.991"
$ws.Range("E42").Value = "  -0.94%  "
$ws.Range("E43").Value = "  -1.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "65.46"
$ws.Range("E44").Value = "  -2.07%  "
$ws.Range("D46").Value = "1.757.54"
$ws.Range("E46").Value = "  -1.18%  "
$ws.Range("E47").Value = "  -4.09%  "
$ws.Range("E48").Value = "  -0.44%  "
$ws.Range("E49").Value = "  +1.07%  "
$ws.Range("E50").Value = "  -0.88%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.57"
$ws.Range("E51").Value = "  +0.88%  "
